$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) figures.
# For Price (column D) cells whose new value parses purely as a
# number (e.g. "1.00", "584.91"), force text format first so Excel
# keeps the exact string instead of silently converting it to a
# floating point number (which would lose trailing zeros / introduce
# rounding artifacts).

$ws.Range("D2").Value = "67.989.78"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "3.336.33"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.91"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.39"
$ws.Range("E6").Value = "  +1.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +1.45%  "
$ws.Range("E9").Value = "  +4.65%  "
$ws.Range("E10").Value = "  +1.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.18"
$ws.Range("E11").Value = "  +6.43%  "
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "695.63"
$ws.Range("E13").Value = "  +4.90%  "
$ws.Range("D14").Value = "3.877.59"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.45"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").Value = "68.006.94"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").Value = "3.341.11"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.52"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("E20").Value = "  +2.81%  "
$ws.Range("E21").Value = "  +1.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.41"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.97"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.24"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.92"
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("E27").Value = "  +2.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.14"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.54"
$ws.Range("E29").Value = "  +1.84%  "
$ws.Range("E30").Value = "  -3.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "575.33"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("E33").Value = "  +2.03%  "
$ws.Range("D34").Value = "3.742.50"
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.45"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.36"
$ws.Range("E37").Value = "  +1.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.31"
$ws.Range("E38").Value = "  +9.24%  "
$ws.Range("E39").Value = "  +3.19%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.17"
$ws.Range("E40").Value = "  +3.18%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.63"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("D42").Value = "0.0₃0676"
$ws.Range("E42").Value = "  +2.21%  "
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.25"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.63"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.11"
$ws.Range("E50").Value = "  +2.89%  "
$ws.Range("E51").Value = "  +0.63%  "

